$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: A1 now carries a "Date" label (was blank) ---
$ws.Range("A1").Value = "Date"

# --- Insert a new row 7 for the "Help" link, pushing everything below down by one ---
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Help"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf")

# --- Column A is a bit wider now to fit the new labels ---
$ws.Columns.Item(1).ColumnWidth = 20.14

# --- OPERATION values renamed throughout the sample rows (now rows 11-26) ---
$ws.Range("A20").Value = "CREATE RELATIONSHIP"
$ws.Range("A21").Value = "INSTALL DEVICE"
$ws.Range("A25").Value = "DELETE RELATIONSHIP"
$ws.Range("A26").Value = "UNINSTALL DEVICE"

# --- Keep the OPERATION dropdown list in sync with the renamed values ---
$ws.Range("A11").Validation.Modify(3, 1, 1, '"CREATE ENTITY,UPDATE ENTITY,DELETE ENTITY,CREATE PROPERTY,UPDATE PROPERTY,DELETE PROPERTY,CREATE RELATIONSHIP,UPDATE RELATIONSHIP,DELETE RELATIONSHIP,INSTALL DEVICE,UNINSTALL DEVICE"')
